# "ran model for jan 23"
# Fill in the "Beat Vegas?" result for the last three games that were
# missing it, then append the model's predictions for the Jan 23 (serial
# 44219) slate of games as new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in previously-missing "Beat Vegas?" values ---
$ws.Cells.Item(108, 7).Value = "Yes"
$ws.Cells.Item(109, 7).Value = "Yes"
$ws.Cells.Item(110, 7).Value = "No"

# --- Append new games for 2021-01-23 (serial date 44219) ---
$dateFmt = "yyyy\-mm\-dd"

$newGames = @(
    @{ Row = 111; Home = "DET"; Away = "PHI"; Spread = 6.5;  Pred = 15.4 },
    @{ Row = 112; Home = "BRK"; Away = "MIA"; Spread = -8;   Pred = -5.8 },
    @{ Row = 113; Home = "MIN"; Away = "NOP"; Spread = 9.5;  Pred = 7 },
    @{ Row = 114; Home = "UTA"; Away = "GSW"; Spread = -7;   Pred = -11.7 },
    @{ Row = 115; Home = "CHI"; Away = "LAL"; Spread = 9.5;  Pred = -6 },
    @{ Row = 116; Home = "DAL"; Away = "HOU"; Spread = -9.5; Pred = 3.6 },
    @{ Row = 117; Home = "PHO"; Away = "DEN"; Spread = 1.5;  Pred = -0.1 }
)

foreach ($game in $newGames) {
    $r = $game.Row
    $ws.Cells.Item($r, 1).Value = 44219
    $ws.Cells.Item($r, 1).NumberFormat = $dateFmt
    $ws.Cells.Item($r, 2).Value = $game.Home
    $ws.Cells.Item($r, 3).Value = $game.Away
    $ws.Cells.Item($r, 4).Value = $game.Spread
    $ws.Cells.Item($r, 5).Value = $game.Pred
    # Spread Difference = Spread - Predicted Spread (matches how the rest of
    # the sheet's "Spread Difference" column is derived).
    $ws.Cells.Item($r, 6).Value = $game.Spread - $game.Pred
}
